$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 11343
$ws.Cells.Item(2, 4).Value = 15576593
$ws.Cells.Item(4, 3).Value = 21391
$ws.Cells.Item(4, 4).Value = 27281916
$ws.Cells.Item(6, 3).Value = 1273
$ws.Cells.Item(6, 4).Value = 1601905
$ws.Cells.Item(7, 3).Value = 59831
$ws.Cells.Item(7, 4).Value = 86032658
$ws.Cells.Item(8, 3).Value = 78160
$ws.Cells.Item(8, 4).Value = 103929819
$ws.Cells.Item(9, 3).Value = 25764
$ws.Cells.Item(9, 4).Value = 35868445
$ws.Cells.Item(10, 3).Value = 64778
$ws.Cells.Item(10, 4).Value = 92766007
$ws.Cells.Item(11, 3).Value = 8951
$ws.Cells.Item(11, 4).Value = 11694394
$ws.Cells.Item(12, 3).Value = 3803
$ws.Cells.Item(12, 4).Value = 5357187
$ws.Cells.Item(13, 3).Value = 14607
$ws.Cells.Item(13, 4).Value = 19772812
$ws.Cells.Item(14, 3).Value = 45897
$ws.Cells.Item(14, 4).Value = 61825461
$ws.Cells.Item(15, 3).Value = 21906
$ws.Cells.Item(15, 4).Value = 28401765
$ws.Cells.Item(17, 3).Value = 42348
$ws.Cells.Item(17, 4).Value = 53635890
$ws.Cells.Item(18, 3).Value = 47654
$ws.Cells.Item(18, 4).Value = 63658204
$ws.Cells.Item(19, 3).Value = 28078
$ws.Cells.Item(19, 4).Value = 34429513
$ws.Cells.Item(20, 3).Value = 49401
$ws.Cells.Item(20, 4).Value = 60043298
$ws.Cells.Item(21, 3).Value = 3823
$ws.Cells.Item(21, 4).Value = 5259827
$ws.Cells.Item(23, 3).Value = 5940
$ws.Cells.Item(23, 4).Value = 7474170
$ws.Cells.Item(25, 3).Value = 497
$ws.Cells.Item(25, 4).Value = 628524
$ws.Cells.Item(26, 3).Value = 14535
$ws.Cells.Item(26, 4).Value = 20820551
$ws.Cells.Item(27, 3).Value = 23832
$ws.Cells.Item(27, 4).Value = 31207092
$ws.Cells.Item(28, 3).Value = 3145
$ws.Cells.Item(28, 4).Value = 4262925
$ws.Cells.Item(29, 3).Value = 18745
$ws.Cells.Item(29, 4).Value = 26636933
$ws.Cells.Item(30, 3).Value = 1649
$ws.Cells.Item(30, 4).Value = 2079252
$ws.Cells.Item(31, 3).Value = 886
$ws.Cells.Item(31, 4).Value = 1210161
$ws.Cells.Item(32, 3).Value = 3199
$ws.Cells.Item(32, 4).Value = 4227980
$ws.Cells.Item(33, 3).Value = 8815
$ws.Cells.Item(33, 4).Value = 11826397
$ws.Cells.Item(34, 3).Value = 4397
$ws.Cells.Item(34, 4).Value = 5586239
$ws.Cells.Item(36, 3).Value = 6377
$ws.Cells.Item(36, 4).Value = 7602595
$ws.Cells.Item(37, 3).Value = 10471
$ws.Cells.Item(37, 4).Value = 13867177
$ws.Cells.Item(38, 3).Value = 7168
$ws.Cells.Item(38, 4).Value = 8565709
$ws.Cells.Item(39, 3).Value = 14961
$ws.Cells.Item(39, 4).Value = 18329676
$ws.Cells.Item(40, 3).Value = 3115
$ws.Cells.Item(40, 4).Value = 4255017
$ws.Cells.Item(42, 3).Value = 7925
$ws.Cells.Item(42, 4).Value = 9692956
$ws.Cells.Item(45, 3).Value = 16032
$ws.Cells.Item(45, 4).Value = 22862440
$ws.Cells.Item(46, 3).Value = 25416
$ws.Cells.Item(46, 4).Value = 33890350
$ws.Cells.Item(47, 3).Value = 4135
$ws.Cells.Item(47, 4).Value = 5701504
$ws.Cells.Item(48, 3).Value = 25117
$ws.Cells.Item(48, 4).Value = 35950614
$ws.Cells.Item(49, 3).Value = 2695
$ws.Cells.Item(49, 4).Value = 3400791
$ws.Cells.Item(50, 3).Value = 1143
$ws.Cells.Item(50, 4).Value = 1611315
$ws.Cells.Item(51, 3).Value = 4761
$ws.Cells.Item(51, 4).Value = 6211030
$ws.Cells.Item(52, 3).Value = 13120
$ws.Cells.Item(52, 4).Value = 17583534
$ws.Cells.Item(53, 3).Value = 5574
$ws.Cells.Item(53, 4).Value = 6933730
$ws.Cells.Item(55, 3).Value = 7111
$ws.Cells.Item(55, 4).Value = 8636631
$ws.Cells.Item(56, 3).Value = 17352
$ws.Cells.Item(56, 4).Value = 23329327
$ws.Cells.Item(57, 3).Value = 8832
$ws.Cells.Item(57, 4).Value = 10560822
$ws.Cells.Item(58, 3).Value = 16713
$ws.Cells.Item(58, 4).Value = 20585537
$ws.Cells.Item(59, 3).Value = 2741
$ws.Cells.Item(59, 4).Value = 3764081
$ws.Cells.Item(60, 3).Value = 4985
$ws.Cells.Item(60, 4).Value = 6318512
$ws.Cells.Item(63, 3).Value = 14316
$ws.Cells.Item(63, 4).Value = 20443973
$ws.Cells.Item(64, 3).Value = 19399
$ws.Cells.Item(64, 4).Value = 25258230
$ws.Cells.Item(65, 3).Value = 4773
$ws.Cells.Item(65, 4).Value = 6672879
$ws.Cells.Item(66, 3).Value = 15175
$ws.Cells.Item(66, 4).Value = 21713608
$ws.Cells.Item(67, 3).Value = 1829
$ws.Cells.Item(67, 4).Value = 2381189
$ws.Cells.Item(68, 3).Value = 760
$ws.Cells.Item(68, 4).Value = 1058572
$ws.Cells.Item(69, 3).Value = 3260
$ws.Cells.Item(69, 4).Value = 4367099
$ws.Cells.Item(70, 3).Value = 8273
$ws.Cells.Item(70, 4).Value = 11155448
$ws.Cells.Item(71, 3).Value = 4501
$ws.Cells.Item(71, 4).Value = 5710077
$ws.Cells.Item(73, 3).Value = 5180
$ws.Cells.Item(73, 4).Value = 6412487
$ws.Cells.Item(74, 3).Value = 8906
$ws.Cells.Item(74, 4).Value = 11747413
$ws.Cells.Item(75, 3).Value = 6359
$ws.Cells.Item(75, 4).Value = 7771141
$ws.Cells.Item(76, 3).Value = 13422
$ws.Cells.Item(76, 4).Value = 16437951
$ws.Cells.Item(77, 3).Value = 2535
$ws.Cells.Item(77, 4).Value = 3503746
$ws.Cells.Item(78, 3).Value = 1991
$ws.Cells.Item(78, 4).Value = 2632024
$ws.Cells.Item(80, 3).Value = 4807
$ws.Cells.Item(80, 4).Value = 6816959
$ws.Cells.Item(81, 3).Value = 5243
$ws.Cells.Item(81, 4).Value = 7300402
$ws.Cells.Item(82, 3).Value = 897
$ws.Cells.Item(82, 4).Value = 1281857
$ws.Cells.Item(83, 3).Value = 5970
$ws.Cells.Item(83, 4).Value = 8575352
$ws.Cells.Item(84, 3).Value = 364
$ws.Cells.Item(84, 4).Value = 481529
$ws.Cells.Item(86, 3).Value = 1413
$ws.Cells.Item(86, 4).Value = 1933979
$ws.Cells.Item(87, 3).Value = 3492
$ws.Cells.Item(87, 4).Value = 4864734
$ws.Cells.Item(88, 3).Value = 2004
$ws.Cells.Item(88, 4).Value = 2564720
$ws.Cells.Item(89, 3).Value = 1024
$ws.Cells.Item(89, 4).Value = 1253231
$ws.Cells.Item(91, 3).Value = 1194
$ws.Cells.Item(91, 4).Value = 1561055
$ws.Cells.Item(93, 3).Value = 5241
$ws.Cells.Item(93, 4).Value = 7295689
$ws.Cells.Item(95, 3).Value = 10394
$ws.Cells.Item(95, 4).Value = 13507582
$ws.Cells.Item(97, 3).Value = 1269
$ws.Cells.Item(97, 4).Value = 1568391
$ws.Cells.Item(98, 3).Value = 30113
$ws.Cells.Item(98, 4).Value = 43252683
$ws.Cells.Item(99, 3).Value = 44012
$ws.Cells.Item(99, 4).Value = 58291029
$ws.Cells.Item(100, 3).Value = 8635
$ws.Cells.Item(100, 4).Value = 11800376
$ws.Cells.Item(101, 3).Value = 30750
$ws.Cells.Item(101, 4).Value = 44214114
$ws.Cells.Item(102, 3).Value = 3641
$ws.Cells.Item(102, 4).Value = 4691827
$ws.Cells.Item(103, 3).Value = 1899
$ws.Cells.Item(103, 4).Value = 2626192
$ws.Cells.Item(104, 3).Value = 5683
$ws.Cells.Item(104, 4).Value = 7622569
$ws.Cells.Item(105, 3).Value = 19521
$ws.Cells.Item(105, 4).Value = 26145707
$ws.Cells.Item(106, 3).Value = 8593
$ws.Cells.Item(106, 4).Value = 10995188
$ws.Cells.Item(108, 3).Value = 10767
$ws.Cells.Item(108, 4).Value = 13057074
$ws.Cells.Item(109, 3).Value = 21551
$ws.Cells.Item(109, 4).Value = 29308098
$ws.Cells.Item(110, 3).Value = 11898
$ws.Cells.Item(110, 4).Value = 14190915
$ws.Cells.Item(111, 3).Value = 30561
$ws.Cells.Item(111, 4).Value = 36663483
$ws.Cells.Item(113, 3).Value = 7145
$ws.Cells.Item(113, 4).Value = 9770238
$ws.Cells.Item(115, 3).Value = 3106
$ws.Cells.Item(115, 4).Value = 4353678
$ws.Cells.Item(117, 3).Value = 63
$ws.Cells.Item(117, 4).Value = 92000
$ws.Cells.Item(118, 3).Value = 5392
$ws.Cells.Item(118, 4).Value = 7792447
$ws.Cells.Item(119, 3).Value = 8500
$ws.Cells.Item(119, 4).Value = 11793048
$ws.Cells.Item(120, 3).Value = 2074
$ws.Cells.Item(120, 4).Value = 2907053
$ws.Cells.Item(121, 3).Value = 5936
$ws.Cells.Item(121, 4).Value = 8489851
$ws.Cells.Item(122, 3).Value = 733
$ws.Cells.Item(122, 4).Value = 1010209
$ws.Cells.Item(124, 3).Value = 1154
$ws.Cells.Item(124, 4).Value = 1595418
$ws.Cells.Item(125, 3).Value = 3079
$ws.Cells.Item(125, 4).Value = 4316735
$ws.Cells.Item(126, 3).Value = 3767
$ws.Cells.Item(126, 4).Value = 5168159
$ws.Cells.Item(127, 3).Value = 1793
$ws.Cells.Item(127, 4).Value = 2289593
$ws.Cells.Item(129, 3).Value = 1308
$ws.Cells.Item(129, 4).Value = 1768917
$ws.Cells.Item(130, 3).Value = 2976
$ws.Cells.Item(130, 4).Value = 3883947
$ws.Cells.Item(132, 3).Value = 2148
$ws.Cells.Item(132, 4).Value = 3015221
$ws.Cells.Item(133, 3).Value = 45
$ws.Cells.Item(133, 4).Value = 71700
$ws.Cells.Item(134, 3).Value = 1010
$ws.Cells.Item(134, 4).Value = 1499599
$ws.Cells.Item(136, 3).Value = 2641
$ws.Cells.Item(136, 4).Value = 4096117
$ws.Cells.Item(137, 3).Value = 2406
$ws.Cells.Item(137, 4).Value = 3584517
$ws.Cells.Item(138, 3).Value = 814
$ws.Cells.Item(138, 4).Value = 1277457
$ws.Cells.Item(139, 3).Value = 1715
$ws.Cells.Item(139, 4).Value = 2637389
$ws.Cells.Item(140, 3).Value = 206
$ws.Cells.Item(140, 4).Value = 312477
$ws.Cells.Item(141, 3).Value = 126
$ws.Cells.Item(141, 4).Value = 195780
$ws.Cells.Item(142, 3).Value = 207
$ws.Cells.Item(142, 4).Value = 307989
$ws.Cells.Item(143, 3).Value = 1063
$ws.Cells.Item(143, 4).Value = 1606506
$ws.Cells.Item(144, 3).Value = 1089
$ws.Cells.Item(144, 4).Value = 1659599
$ws.Cells.Item(145, 3).Value = 461
$ws.Cells.Item(145, 4).Value = 680011
$ws.Cells.Item(146, 3).Value = 518
$ws.Cells.Item(146, 4).Value = 756563
$ws.Cells.Item(147, 3).Value = 400
$ws.Cells.Item(147, 4).Value = 586656
$ws.Cells.Item(148, 3).Value = 813
$ws.Cells.Item(148, 4).Value = 1171025
$ws.Cells.Item(149, 3).Value = 2849
$ws.Cells.Item(149, 4).Value = 3881160
$ws.Cells.Item(151, 3).Value = 7954
$ws.Cells.Item(151, 4).Value = 9833312
$ws.Cells.Item(154, 3).Value = 24826
$ws.Cells.Item(154, 4).Value = 35271861
$ws.Cells.Item(155, 3).Value = 42096
$ws.Cells.Item(155, 4).Value = 53474708
$ws.Cells.Item(156, 3).Value = 13434
$ws.Cells.Item(156, 4).Value = 18497667
$ws.Cells.Item(157, 3).Value = 29451
$ws.Cells.Item(157, 4).Value = 42366738
$ws.Cells.Item(158, 3).Value = 3672
$ws.Cells.Item(158, 4).Value = 4752005
$ws.Cells.Item(159, 3).Value = 1897
$ws.Cells.Item(159, 4).Value = 2658574
$ws.Cells.Item(160, 3).Value = 5157
$ws.Cells.Item(160, 4).Value = 6910892
$ws.Cells.Item(161, 3).Value = 18692
$ws.Cells.Item(161, 4).Value = 25348196
$ws.Cells.Item(162, 3).Value = 8455
$ws.Cells.Item(162, 4).Value = 10626897
$ws.Cells.Item(164, 3).Value = 9469
$ws.Cells.Item(164, 4).Value = 11685032
$ws.Cells.Item(165, 3).Value = 22115
$ws.Cells.Item(165, 4).Value = 29945130
$ws.Cells.Item(166, 3).Value = 11677
$ws.Cells.Item(166, 4).Value = 14162568
$ws.Cells.Item(167, 3).Value = 29004
$ws.Cells.Item(167, 4).Value = 34002630
$ws.Cells.Item(168, 3).Value = 876
$ws.Cells.Item(168, 4).Value = 1211686
$ws.Cells.Item(170, 3).Value = 21837
$ws.Cells.Item(170, 4).Value = 29110994
$ws.Cells.Item(171, 3).Value = 398
$ws.Cells.Item(171, 4).Value = 581356
$ws.Cells.Item(172, 3).Value = 997
$ws.Cells.Item(172, 4).Value = 1363889
$ws.Cells.Item(173, 3).Value = 65972
$ws.Cells.Item(173, 4).Value = 95137767
$ws.Cells.Item(174, 3).Value = 118551
$ws.Cells.Item(174, 4).Value = 162282381
$ws.Cells.Item(175, 3).Value = 183057
$ws.Cells.Item(175, 4).Value = 262310553
$ws.Cells.Item(176, 3).Value = 78932
$ws.Cells.Item(176, 4).Value = 115933244
$ws.Cells.Item(177, 3).Value = 37235
$ws.Cells.Item(177, 4).Value = 51011132
$ws.Cells.Item(178, 3).Value = 8432
$ws.Cells.Item(178, 4).Value = 12062927
$ws.Cells.Item(179, 3).Value = 21189
$ws.Cells.Item(179, 4).Value = 29955998
$ws.Cells.Item(180, 3).Value = 137027
$ws.Cells.Item(180, 4).Value = 188209055
$ws.Cells.Item(181, 3).Value = 37712
$ws.Cells.Item(181, 4).Value = 50711199
$ws.Cells.Item(183, 3).Value = 37973
$ws.Cells.Item(183, 4).Value = 46696135
$ws.Cells.Item(184, 3).Value = 58416
$ws.Cells.Item(184, 4).Value = 78811452
$ws.Cells.Item(185, 3).Value = 60597
$ws.Cells.Item(185, 4).Value = 77805701
$ws.Cells.Item(186, 3).Value = 62132
$ws.Cells.Item(186, 4).Value = 80899739
$ws.Cells.Item(187, 3).Value = 4863
$ws.Cells.Item(187, 4).Value = 6448946
$ws.Cells.Item(189, 3).Value = 4078
$ws.Cells.Item(189, 4).Value = 5446737
$ws.Cells.Item(192, 3).Value = 7556
$ws.Cells.Item(192, 4).Value = 10976240
$ws.Cells.Item(193, 3).Value = 13279
$ws.Cells.Item(193, 4).Value = 18023383
$ws.Cells.Item(194, 3).Value = 1851
$ws.Cells.Item(194, 4).Value = 2591636
$ws.Cells.Item(195, 3).Value = 7128
$ws.Cells.Item(195, 4).Value = 10089961
$ws.Cells.Item(196, 3).Value = 959
$ws.Cells.Item(196, 4).Value = 1296049
$ws.Cells.Item(197, 3).Value = 398
$ws.Cells.Item(197, 4).Value = 574998
$ws.Cells.Item(198, 3).Value = 1615
$ws.Cells.Item(198, 4).Value = 2242183
$ws.Cells.Item(199, 3).Value = 4212
$ws.Cells.Item(199, 4).Value = 5912647
$ws.Cells.Item(200, 3).Value = 2569
$ws.Cells.Item(200, 4).Value = 3538386
$ws.Cells.Item(201, 3).Value = 3419
$ws.Cells.Item(201, 4).Value = 4450152
$ws.Cells.Item(202, 3).Value = 5196
$ws.Cells.Item(202, 4).Value = 7343738
$ws.Cells.Item(203, 3).Value = 2394
$ws.Cells.Item(203, 4).Value = 3149447
$ws.Cells.Item(204, 3).Value = 5041
$ws.Cells.Item(204, 4).Value = 6416267
$ws.Cells.Item(205, 3).Value = 1805
$ws.Cells.Item(205, 4).Value = 2279392
$ws.Cells.Item(206, 3).Value = 1989
$ws.Cells.Item(206, 4).Value = 2673625
$ws.Cells.Item(209, 3).Value = 3151
$ws.Cells.Item(209, 4).Value = 4550251
$ws.Cells.Item(210, 3).Value = 5359
$ws.Cells.Item(210, 4).Value = 7344178
$ws.Cells.Item(211, 3).Value = 1809
$ws.Cells.Item(211, 4).Value = 2549647
$ws.Cells.Item(212, 3).Value = 3177
$ws.Cells.Item(212, 4).Value = 4552203
$ws.Cells.Item(213, 3).Value = 526
$ws.Cells.Item(213, 4).Value = 708844
$ws.Cells.Item(215, 3).Value = 672
$ws.Cells.Item(215, 4).Value = 942736
$ws.Cells.Item(216, 3).Value = 2485
$ws.Cells.Item(216, 4).Value = 3437203
$ws.Cells.Item(217, 3).Value = 2416
$ws.Cells.Item(217, 4).Value = 3300000
$ws.Cells.Item(218, 3).Value = 1308
$ws.Cells.Item(218, 4).Value = 1696287
$ws.Cells.Item(220, 3).Value = 1000
$ws.Cells.Item(220, 4).Value = 1351157
$ws.Cells.Item(221, 3).Value = 2735
$ws.Cells.Item(221, 4).Value = 3576047
$ws.Cells.Item(223, 3).Value = 3074
$ws.Cells.Item(223, 4).Value = 4757612
$ws.Cells.Item(224, 3).Value = 616
$ws.Cells.Item(224, 4).Value = 957166
$ws.Cells.Item(226, 3).Value = 2000
$ws.Cells.Item(226, 4).Value = 3162797
$ws.Cells.Item(227, 3).Value = 9164
$ws.Cells.Item(227, 4).Value = 14193079
$ws.Cells.Item(228, 3).Value = 1768
$ws.Cells.Item(228, 4).Value = 2779321
$ws.Cells.Item(229, 3).Value = 1024
$ws.Cells.Item(229, 4).Value = 1615245
$ws.Cells.Item(230, 3).Value = 117
$ws.Cells.Item(230, 4).Value = 180476
$ws.Cells.Item(231, 3).Value = 29
$ws.Cells.Item(231, 4).Value = 44600
$ws.Cells.Item(233, 3).Value = 363
$ws.Cells.Item(233, 4).Value = 556882
$ws.Cells.Item(234, 3).Value = 353
$ws.Cells.Item(234, 4).Value = 558899
$ws.Cells.Item(235, 3).Value = 257
$ws.Cells.Item(235, 4).Value = 417794
$ws.Cells.Item(236, 3).Value = 210
$ws.Cells.Item(236, 4).Value = 320430
$ws.Cells.Item(237, 3).Value = 168
$ws.Cells.Item(237, 4).Value = 260245
$ws.Cells.Item(238, 3).Value = 344
$ws.Cells.Item(238, 4).Value = 525365
$ws.Cells.Item(239, 3).Value = 3213
$ws.Cells.Item(239, 4).Value = 4383583
$ws.Cells.Item(241, 3).Value = 5538
$ws.Cells.Item(241, 4).Value = 7007587
$ws.Cells.Item(244, 3).Value = 15425
$ws.Cells.Item(244, 4).Value = 22056944
$ws.Cells.Item(245, 3).Value = 26496
$ws.Cells.Item(245, 4).Value = 34521726
$ws.Cells.Item(246, 3).Value = 4906
$ws.Cells.Item(246, 4).Value = 6764312
$ws.Cells.Item(247, 3).Value = 20301
$ws.Cells.Item(247, 4).Value = 29095723
$ws.Cells.Item(248, 3).Value = 1958
$ws.Cells.Item(248, 4).Value = 2476557
$ws.Cells.Item(249, 3).Value = 1176
$ws.Cells.Item(249, 4).Value = 1634327
$ws.Cells.Item(250, 3).Value = 3828
$ws.Cells.Item(250, 4).Value = 5051000
$ws.Cells.Item(251, 3).Value = 11047
$ws.Cells.Item(251, 4).Value = 14932026
$ws.Cells.Item(252, 3).Value = 5123
$ws.Cells.Item(252, 4).Value = 6470532
$ws.Cells.Item(254, 3).Value = 6125
$ws.Cells.Item(254, 4).Value = 7323591
$ws.Cells.Item(255, 3).Value = 10582
$ws.Cells.Item(255, 4).Value = 13962216
$ws.Cells.Item(256, 3).Value = 7725
$ws.Cells.Item(256, 4).Value = 9412717
$ws.Cells.Item(257, 3).Value = 17704
$ws.Cells.Item(257, 4).Value = 21492143
$ws.Cells.Item(258, 3).Value = 12034
$ws.Cells.Item(258, 4).Value = 16636627
$ws.Cells.Item(260, 3).Value = 17011
$ws.Cells.Item(260, 4).Value = 20904082
$ws.Cells.Item(262, 3).Value = 927
$ws.Cells.Item(262, 4).Value = 1052669
$ws.Cells.Item(263, 3).Value = 44747
$ws.Cells.Item(263, 4).Value = 63464258
$ws.Cells.Item(264, 3).Value = 64915
$ws.Cells.Item(264, 4).Value = 85143058
$ws.Cells.Item(265, 3).Value = 12374
$ws.Cells.Item(265, 4).Value = 17012590
$ws.Cells.Item(266, 3).Value = 45154
$ws.Cells.Item(266, 4).Value = 63780877
$ws.Cells.Item(267, 3).Value = 5686
$ws.Cells.Item(267, 4).Value = 7348972
$ws.Cells.Item(268, 3).Value = 2893
$ws.Cells.Item(268, 4).Value = 4050165
$ws.Cells.Item(269, 3).Value = 11193
$ws.Cells.Item(269, 4).Value = 14771135
$ws.Cells.Item(270, 3).Value = 32109
$ws.Cells.Item(270, 4).Value = 43075292
$ws.Cells.Item(271, 3).Value = 16272
$ws.Cells.Item(271, 4).Value = 20435668
$ws.Cells.Item(273, 3).Value = 17840
$ws.Cells.Item(273, 4).Value = 20820193
$ws.Cells.Item(274, 3).Value = 34370
$ws.Cells.Item(274, 4).Value = 45315386
$ws.Cells.Item(275, 3).Value = 19065
$ws.Cells.Item(275, 4).Value = 22912385
$ws.Cells.Item(276, 3).Value = 40751
$ws.Cells.Item(276, 4).Value = 48705822
$ws.Cells.Item(277, 3).Value = 12971
$ws.Cells.Item(277, 4).Value = 17459195
$ws.Cells.Item(279, 3).Value = 18547
$ws.Cells.Item(279, 4).Value = 22938384
$ws.Cells.Item(282, 3).Value = 57504
$ws.Cells.Item(282, 4).Value = 81448488
$ws.Cells.Item(283, 3).Value = 73085
$ws.Cells.Item(283, 4).Value = 95787446
$ws.Cells.Item(284, 3).Value = 14239
$ws.Cells.Item(284, 4).Value = 19233790
$ws.Cells.Item(285, 3).Value = 54238
$ws.Cells.Item(285, 4).Value = 76796394
$ws.Cells.Item(286, 3).Value = 6886
$ws.Cells.Item(286, 4).Value = 8787272
$ws.Cells.Item(287, 3).Value = 3026
$ws.Cells.Item(287, 4).Value = 4183822
$ws.Cells.Item(288, 3).Value = 12699
$ws.Cells.Item(288, 4).Value = 16939725
$ws.Cells.Item(289, 3).Value = 36877
$ws.Cells.Item(289, 4).Value = 49812067
$ws.Cells.Item(290, 3).Value = 17887
$ws.Cells.Item(290, 4).Value = 22537252
$ws.Cells.Item(292, 3).Value = 21863
$ws.Cells.Item(292, 4).Value = 25649209
$ws.Cells.Item(293, 3).Value = 37621
$ws.Cells.Item(293, 4).Value = 49819144
$ws.Cells.Item(294, 3).Value = 22222
$ws.Cells.Item(294, 4).Value = 26680764
$ws.Cells.Item(295, 3).Value = 41886
$ws.Cells.Item(295, 4).Value = 48655023
$ws.Cells.Item(296, 3).Value = 3984
$ws.Cells.Item(296, 4).Value = 5565510
$ws.Cells.Item(298, 3).Value = 7598
$ws.Cells.Item(298, 4).Value = 9392989
$ws.Cells.Item(300, 3).Value = 548
$ws.Cells.Item(300, 4).Value = 698827
$ws.Cells.Item(301, 3).Value = 18046
$ws.Cells.Item(301, 4).Value = 25853783
$ws.Cells.Item(302, 3).Value = 28263
$ws.Cells.Item(302, 4).Value = 36922644
$ws.Cells.Item(303, 3).Value = 6330
$ws.Cells.Item(303, 4).Value = 8830258
$ws.Cells.Item(304, 3).Value = 21747
$ws.Cells.Item(304, 4).Value = 31150236
$ws.Cells.Item(305, 3).Value = 3185
$ws.Cells.Item(305, 4).Value = 4092988
$ws.Cells.Item(306, 3).Value = 1773
$ws.Cells.Item(306, 4).Value = 2498542
$ws.Cells.Item(307, 3).Value = 5540
$ws.Cells.Item(307, 4).Value = 7399344
$ws.Cells.Item(308, 3).Value = 16822
$ws.Cells.Item(308, 4).Value = 22682353
$ws.Cells.Item(309, 3).Value = 6288
$ws.Cells.Item(309, 4).Value = 8113816
$ws.Cells.Item(310, 3).Value = 8096
$ws.Cells.Item(310, 4).Value = 9692239
$ws.Cells.Item(311, 3).Value = 19164
$ws.Cells.Item(311, 4).Value = 25234080
$ws.Cells.Item(312, 3).Value = 10449
$ws.Cells.Item(312, 4).Value = 12745101
$ws.Cells.Item(313, 3).Value = 20454
$ws.Cells.Item(313, 4).Value = 25116438
$ws.Cells.Item(314, 3).Value = 5371
$ws.Cells.Item(314, 4).Value = 7300311
$ws.Cells.Item(316, 3).Value = 17655
$ws.Cells.Item(316, 4).Value = 22835177
$ws.Cells.Item(317, 3).Value = 77
$ws.Cells.Item(317, 4).Value = 101863
$ws.Cells.Item(318, 3).Value = 344
$ws.Cells.Item(318, 4).Value = 444542
$ws.Cells.Item(319, 3).Value = 50433
$ws.Cells.Item(319, 4).Value = 71905265
$ws.Cells.Item(320, 3).Value = 77427
$ws.Cells.Item(320, 4).Value = 102815598
$ws.Cells.Item(321, 3).Value = 23675
$ws.Cells.Item(321, 4).Value = 33200750
$ws.Cells.Item(322, 3).Value = 51715
$ws.Cells.Item(322, 4).Value = 74508868
$ws.Cells.Item(323, 3).Value = 7465
$ws.Cells.Item(323, 4).Value = 9741088
$ws.Cells.Item(324, 3).Value = 3394
$ws.Cells.Item(324, 4).Value = 4776887
$ws.Cells.Item(325, 3).Value = 15730
$ws.Cells.Item(325, 4).Value = 21666197
$ws.Cells.Item(326, 3).Value = 39708
$ws.Cells.Item(326, 4).Value = 53863388
$ws.Cells.Item(327, 3).Value = 21672
$ws.Cells.Item(327, 4).Value = 28090726
$ws.Cells.Item(329, 3).Value = 23185
$ws.Cells.Item(329, 4).Value = 28060651
$ws.Cells.Item(330, 3).Value = 34573
$ws.Cells.Item(330, 4).Value = 45942459
$ws.Cells.Item(331, 3).Value = 20802
$ws.Cells.Item(331, 4).Value = 25864853
$ws.Cells.Item(332, 3).Value = 41737
$ws.Cells.Item(332, 4).Value = 49757825
